$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 04:44"

# Helper: write a full country data row (A..H)
function Set-CountryRow($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Peru overtakes Sudafrica (rows 8-9) ---
Set-CountryRow 8 "Peru" 621997 0 429662 164058 0 0 28277
Set-CountryRow 9 "Sudafrica" 618286 0 531338 73320 0 0 13628

# --- Bolivia data refresh (row 29), no reordering ---
Set-CountryRow 29 "Bolivia" 113129 1035 52521 55817 0 65 4791

# --- Belgica overtakes Rumania (rows 40-41) ---
Set-CountryRow 40 "Belgica" 83500 470 18360 55256 0 5 9884
Set-CountryRow 41 "Rumania" 83150 0 36677 43014 0 0 3459

# --- Honduras jumps ahead of Marruecos & Portugal (rows 49-51) ---
Set-CountryRow 49 "Honduras" 57669 1020 9586 46280 0 56 1803
Set-CountryRow 50 "Marruecos" 57085 0 41901 14173 0 0 1011
Set-CountryRow 51 "Portugal" 56673 0 41357 13507 0 0 1809

# --- Australia data refresh (row 72), no reordering ---
Set-CountryRow 72 "Australia" 25446 124 20367 4496 0 11 583

# --- Corea del Sur overtakes Bosnia y Herzegovina (rows 77-78) ---
Set-CountryRow 77 "Corea del Sur" 19077 371 14551 4210 0 3 316
Set-CountryRow 78 "Bosnia y Herzegovina" 18920 0 12480 5858 0 0 582

# --- Paraguay overtakes Libano (rows 82-83) ---
Set-CountryRow 82 "Paraguay" 15290 0 8348 6677 0 0 265
Set-CountryRow 83 "Libano" 14937 0 4133 10658 0 0 146

# --- Camboya minor correction (row 184) ---
$ws.Cells.Item(184, 4).Value = 265
$ws.Cells.Item(184, 5).Value = 8
